# Apply updated cryptocurrency price/volume data to worksheet cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.111.32'
$ws.Range("D3").Value = '3.477.10'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.52'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.02%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -1.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.71'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.02%  '
$ws.Range("E10").Value = '  -0.92%  '
$ws.Range("E11").Value = '  -0.04%  '
$ws.Range("D12").Value = '4.066.77'
$ws.Range("E12").Value = '  -0.37%  '
$ws.Range("E13").Value = '  +0.06%  '
$ws.Range("E14").Value = '  -2.46%  '
$ws.Range("D15").Value = '3.475.25'
$ws.Range("E15").Value = '  -0.41%  '
$ws.Range("D16").Value = '64.088.61'
$ws.Range("E16").Value = '  -0.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.12'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.63%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '9.96'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("E19").Value = '  -1.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '384.76'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.44%  '
$ws.Range("E22").Value = '  -0.56%  '
$ws.Range("D23").Value = '3.616.78'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.49'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("E25").Value = '  +0.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000112'
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Value = '  +0.28%  '
$ws.Range("E28").Value = '  -0.41%  '
$ws.Range("E29").Value = '  -4.23%  '
$ws.Range("E30").Value = '  -3.22%  '
$ws.Range("E31").Value = '  +2.90%  '
$ws.Range("E32").Value = '  -4.35%  '
$ws.Range("D33").Value = '3.504.49'
$ws.Range("E33").Value = '  -0.21%  '
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("E35").Value = '  -1.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.20'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.13%  '
$ws.Range("E37").Value = '  -2.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '162.79'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.07%  '
$ws.Range("E39").Value = '  -3.64%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0777'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.50%  '
$ws.Range("E41").Value = '  -0.94%  '
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.54'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.79%  '
$ws.Range("E44").Value = '  -1.28%  '
$ws.Range("E45").Value = '  -2.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.53'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.06%  '
$ws.Range("E47").Value = '  -3.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.72'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.897'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.37%  '
$ws.Range("D50").Value = '2.328.70'
$ws.Range("E50").Value = '  -5.24%  '
$ws.Range("E51").Value = '  -2.67%  '
